$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "28.073.65"
Set-TextValue $ws "E2" "  -0.23%  "
Set-TextValue $ws "D3" "1.758.91"
Set-TextValue $ws "E3" "  -1.13%  "
Set-TextValue $ws "D4" "0.9984"
Set-TextValue $ws "E4" "  -0.74%  "
Set-TextValue $ws "D5" "333.49"
Set-TextValue $ws "E5" "  -1.53%  "
Set-TextValue $ws "D6" "0.9948"
Set-TextValue $ws "E6" "  -0.58%  "
Set-TextValue $ws "D7" "0.3876"
Set-TextValue $ws "E7" "  +1.76%  "
Set-TextValue $ws "D8" "0.3398"
Set-TextValue $ws "E8" "  -1.17%  "
Set-TextValue $ws "D9" "45.51"
Set-TextValue $ws "E9" "  -3.28%  "
Set-TextValue $ws "D10" "1.123"
Set-TextValue $ws "E10" "  -2.33%  "
Set-TextValue $ws "D11" "0.07210"
Set-TextValue $ws "E11" "  -2.30%  "
Set-TextValue $ws "D12" "0.9954"
Set-TextValue $ws "E12" "  -0.47%  "
Set-TextValue $ws "D13" "22.34"
Set-TextValue $ws "E13" "  -3.46%  "
Set-TextValue $ws "D14" "6.153"
Set-TextValue $ws "E14" "  -4.14%  "
Set-TextValue $ws "D15" "1.747.89"
Set-TextValue $ws "E15" "  -1.82%  "
Set-TextValue $ws "D16" "7.031"
Set-TextValue $ws "E16" "  -3.16%  "
Set-TextValue $ws "D17" "0.00001056"
Set-TextValue $ws "E17" "  -1.46%  "
Set-TextValue $ws "D18" "0.06592"
Set-TextValue $ws "E18" "  -0.85%  "
Set-TextValue $ws "D19" "80.64"
Set-TextValue $ws "E19" "  -1.97%  "
Set-TextValue $ws "D20" "0.9956"
Set-TextValue $ws "E20" "  -0.49%  "
Set-TextValue $ws "D21" "16.97"
Set-TextValue $ws "E21" "  -2.88%  "
Set-TextValue $ws "D22" "6.195"
Set-TextValue $ws "E22" "  -3.83%  "
Set-TextValue $ws "D23" "28.067.49"
Set-TextValue $ws "E23" "  -0.40%  "
Set-TextValue $ws "D24" "11.60"
Set-TextValue $ws "E24" "  -3.54%  "
Set-TextValue $ws "D25" "2.374"
Set-TextValue $ws "E25" "  +0.29%  "
Set-TextValue $ws "D26" "154.17"
Set-TextValue $ws "E26" "  +0.29%  "
Set-TextValue $ws "D27" "19.90"
Set-TextValue $ws "E27" "  -3.99%  "
Set-TextValue $ws "D28" "2.320"
Set-TextValue $ws "E28" "  -3.72%  "
Set-TextValue $ws "D29" "1.948.46"
Set-TextValue $ws "E29" "  -1.73%  "
Set-TextValue $ws "D30" "1.288"
Set-TextValue $ws "E30" "  -10.28%  "
Set-TextValue $ws "D31" "128.87"
Set-TextValue $ws "E31" "  -5.59%  "
Set-TextValue $ws "D32" "4.062"
Set-TextValue $ws "E32" "  +3.38%  "
Set-TextValue $ws "D33" "5.847"
Set-TextValue $ws "E33" "  -4.36%  "
Set-TextValue $ws "D34" "0.08672"
Set-TextValue $ws "E34" "  -2.19%  "
Set-TextValue $ws "D35" "12.09"
Set-TextValue $ws "E35" "  -5.16%  "
Set-TextValue $ws "D36" "5.143"
Set-TextValue $ws "E36" "  -2.89%  "
Set-TextValue $ws "D37" "0.02283"
Set-TextValue $ws "E37" "  -5.89%  "
Set-TextValue $ws "D38" "0.06151"
Set-TextValue $ws "E38" "  -3.03%  "
Set-TextValue $ws "D39" "0.6476"
Set-TextValue $ws "E39" "  -5.12%  "
Set-TextValue $ws "D40" "1.508"
Set-TextValue $ws "E40" "  +0.64%  "
Set-TextValue $ws "D41" "0.2104"
Set-TextValue $ws "E41" "  -2.79%  "
Set-TextValue $ws "D42" "1.200"
Set-TextValue $ws "E42" "  -3.14%  "
Set-TextValue $ws "D43" "0.9964"
Set-TextValue $ws "E43" "  -0.39%  "
Set-TextValue $ws "D44" "7.850"
Set-TextValue $ws "E44" "  -5.04%  "
Set-TextValue $ws "D45" "13.66"
Set-TextValue $ws "E45" "  -3.71%  "
Set-TextValue $ws "D46" "3.812"
Set-TextValue $ws "E46" "  -1.63%  "
Set-TextValue $ws "D47" "0.5997"
Set-TextValue $ws "E47" "  -4.28%  "
Set-TextValue $ws "D48" "126.29"
Set-TextValue $ws "E48" "  -4.79%  "
Set-TextValue $ws "D49" "1.979"
Set-TextValue $ws "E49" "  -5.03%  "
Set-TextValue $ws "B50" "EOS"
Set-TextValue $ws "C50" "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextValue $ws "D50" "1.154"
Set-TextValue $ws "E50" "  -4.11%  "
Set-TextValue $ws "B51" "Cronos"
Set-TextValue $ws "C51" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D51" "0.06996"
Set-TextValue $ws "E51" "  -6.03%  "

Write-Output "Update complete"